$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First and second models have finished training -> mark their Status as "done"
$ws.Range("D2").Value = "done"
$ws.Range("D3").Value = "done"

# Third model's status remains "trainable" (unchanged)

# Reflect that the last selected cell before saving was D4
$ws.Range("D4").Select()
